$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "58.190.50"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "2.500.69"
$ws.Range("E3").Value = "  +1.46%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "519.37"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").Value = "131.73"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").Value = "2.520.61"
$ws.Range("E9").Value = "  +2.28%  "

$ws.Range("D10").Value = "0.0972"
$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  -2.39%  "

$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").Value = "2.944.09"
$ws.Range("E14").Value = "  +1.47%  "

$ws.Range("D15").Value = "58.141.63"
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").Value = "2.511.00"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").Value = "324.08"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "6.05"
$ws.Range("E22").Value = "  +5.70%  "

$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").Value = "63.51"
$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").Value = "0.161"
$ws.Range("E26").Value = "  +1.27%  "

$ws.Range("D27").Value = "0.991"
$ws.Range("E27").Value = "  -0.91%  "

$ws.Range("D28").Value = "7.35"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").Value = "168.32"
$ws.Range("E30").Value = "  +1.40%  "

$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  +2.69%  "

$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").Value = "0.995"
$ws.Range("E35").Value = "  -0.32%  "

$ws.Range("D36").Value = "18.07"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("E37").Value = "  -2.59%  "

$ws.Range("E38").Value = "  -0.50%  "

$ws.Range("D39").Value = "36.84"
$ws.Range("E39").Value = "  +0.84%  "

$ws.Range("E40").Value = "  -0.60%  "

$ws.Range("D41").Value = "0.776"
$ws.Range("E41").Value = "  -1.45%  "

$ws.Range("D42").Value = "280.34"
$ws.Range("E42").Value = "  +3.33%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "5.13"
$ws.Range("E43").Value = "  +2.92%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.44"
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("E45").Value = "  +1.64%  "

$ws.Range("D46").Value = "122.62"
$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("D47").Value = "0.0920"
$ws.Range("E47").Value = "  +1.86%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "0.0497"
$ws.Range("E48").Value = "  +2.05%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "17.77"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").Value = "17.02"
$ws.Range("E51").Value = "  +0.15%  "
